$d = $word.ActiveDocument

# ===================================================================
# Paragraph 7: "Rilevazione (...)"
#   - ", ClasseRilevCLASSE, IndividuoRilevPERSONA)" -> ", RespRilevRESPONSABILE) "
#   - ", DataIns, RespRil" + "ev" + ", ModAcquisizione, InfoAmb"
#       -> ", DataIns, ModAcquisizione, InfoAmb"
# ===================================================================
$p7 = $d.Paragraphs(7)
$p7Start = $p7.Range.Start
$p7Text = $p7.Range.Text

# --- rightmost edit first: replace the red Classe/Individuo segment ---
$iB = $p7Text.IndexOf(", ClasseRilev")
$iBend = $p7Text.IndexOf(")", $iB) + 1
$absB1 = $p7Start + $iB
$absB2 = $p7Start + $iBend
$rngB = $d.Range($absB1, $absB2)
$rngB.Text = ""

$cursor = $absB1

$ins1 = $d.Range($cursor, $cursor)
$t1 = ", RespRilev"
$ins1.InsertAfter($t1)
$cursor = $cursor + $t1.Length

$ins2 = $d.Range($cursor, $cursor)
$t2 = "RESPONSABILE"
$ins2.InsertAfter($t2)
$rngSuper = $d.Range($cursor, $cursor + $t2.Length)
$rngSuper.Font.Superscript = $true
$cursor = $cursor + $t2.Length

$ins3 = $d.Range($cursor, $cursor)
$t3 = ") "
$ins3.InsertAfter($t3)
$cursor = $cursor + $t3.Length

# --- leftward edit: remove the duplicated ", RespRilev" / "ev" segment ---
$p7Text2 = $p7.Range.Text
$iC = $p7Text2.IndexOf(", DataIns, RespRil")
$iCend = $p7Text2.IndexOf(", ModAcquisizione, InfoAmb", $iC) + ", ModAcquisizione, InfoAmb".Length
$absC1 = $p7Start + $iC
$absC2 = $p7Start + $iCend
$rngC = $d.Range($absC1, $absC2)
$rngC.Text = ", DataIns, ModAcquisizione, InfoAmb"

# ===================================================================
# Paragraph 6: "Responsabile (...)"
#   - "Nominativo" -> "Tipo"
#   - "Cod" + "Resp" (two underlined runs) -> "CodResp" (single underlined run)
# ===================================================================
$p6 = $d.Paragraphs(6)
$p6Start = $p6.Range.Start
$p6Text = $p6.Range.Text

# --- rightmost edit first: Nominativo -> Tipo ---
$iN = $p6Text.IndexOf("Nominativo")
$absN1 = $p6Start + $iN
$absN2 = $absN1 + "Nominativo".Length
$rngN = $d.Range($absN1, $absN2)
$rngN.Text = "Tipo"

# --- leftward edit: merge "Cod" + "Resp" into a single run "CodResp" ---
$p6Text2 = $p6.Range.Text
$iCR = $p6Text2.IndexOf("Cod")
$absCR1 = $p6Start + $iCR
$absCR2 = $absCR1 + "CodResp".Length
$rngCR = $d.Range($absCR1, $absCR2)
$rngCR.Text = ""
$insCR = $d.Range($absCR1, $absCR1)
$insCR.InsertAfter("CodResp")
$rngCR2 = $d.Range($absCR1, $absCR1 + 7)
$rngCR2.Font.Underline = 1

# ===================================================================
# Paragraph 5: "Classe (...)"
#   - ", Ordine, TipoScuola" -> ", Nome, Ordine, TipoScuola"
# ===================================================================
$p5 = $d.Paragraphs(5)
$p5Start = $p5.Range.Start
$p5Text = $p5.Range.Text

$iO = $p5Text.IndexOf(", Ordine, TipoScuola")
$absO1 = $p5Start + $iO + 1
$insO = $d.Range($absO1, $absO1)
$insO.InsertAfter(" Nome,")

# ===================================================================
# Paragraph 4: "Scuola (...)"
#   - remove the space between "Prov," and "CicloIstruz"
# ===================================================================
$p4 = $d.Paragraphs(4)
$p4Start = $p4.Range.Start
$p4Text = $p4.Range.Text

$iP = $p4Text.IndexOf("Prov, CicloIstruz")
$absSpace = $p4Start + $iP + "Prov,".Length
$rngSpace = $d.Range($absSpace, $absSpace + 1)
$rngSpace.Text = ""

Write-Output ("Para4: " + $d.Paragraphs(4).Range.Text)
Write-Output ("Para5: " + $d.Paragraphs(5).Range.Text)
Write-Output ("Para6: " + $d.Paragraphs(6).Range.Text)
Write-Output ("Para7: " + $d.Paragraphs(7).Range.Text)
